$p = $ppt.ActivePresentation

$oldDate = "18/05/2015"
$newDate = "19/05/2015"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master "Date Placeholder" shape
Update-DateShapes $p.SlideMaster.Shapes

# Every slide layout's "Date Placeholder" shape
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateShapes $layouts.Item($j).Shapes
}

# Notes master "Date Placeholder" shape
Update-DateShapes $p.NotesMaster.Shapes
